$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value2 = 364.33334
$ws.Range("I4").Value2 = 127.75
$ws.Range("J4").Value2 = 837.5
$ws.Range("K4").Value2 = 127.75
$ws.Range("L4").Value2 = 837.5
$ws.Range("M4").Value2 = -13.75
$ws.Range("N4").Value2 = -1065.5
$ws.Range("H6").Value2 = 1171.2222
$ws.Range("I6").Value2 = 498.4
$ws.Range("J6").Value2 = 2012.25
$ws.Range("K6").Value2 = 1495.2
$ws.Range("L6").Value2 = 6036.75
$ws.Range("M6").Value2 = -1383.2
$ws.Range("N6").Value2 = -6260.75
$ws.Range("H74").Value2 = 11598.4
$ws.Range("I74").Value2 = 9498
$ws.Range("J74").Value2 = 12123.5
$ws.Range("K74").Value2 = 9498
$ws.Range("L74").Value2 = 12123.5
$ws.Range("M74").Value2 = -8562
$ws.Range("N74").Value2 = -13995.5
$ws.Range("H77").Value2 = 11598.4
$ws.Range("I77").Value2 = 9498
$ws.Range("J77").Value2 = 12123.5
$ws.Range("K77").Value2 = 47490
$ws.Range("L77").Value2 = 60617.5
$ws.Range("M77").Value2 = -42810
$ws.Range("N77").Value2 = -69977.5
$ws.Range("H107").Value2 = 200
$ws.Range("I107").Value2 = 200
$ws.Range("K107").Value2 = 200
$ws.Range("M107").Value2 = 1720
$ws.Range("H132").Value2 = 10620.35
$ws.Range("I132").Value2 = 10620.35
$ws.Range("J132").Value2 = 0
$ws.Range("K132").Value2 = 31861.05
$ws.Range("L132").Value2 = 0
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value2 = -29331.05
$ws.Range("H135").Value2 = 1619.6923
$ws.Range("I135").Value2 = 1172.3
$ws.Range("K135").Value2 = 10550.7
$ws.Range("M135").Value2 = -8015.699999999999
$ws.Range("H138").Value2 = 2465.5557
$ws.Range("J138").Value2 = 3721.8572
$ws.Range("L138").Value2 = 11165.5716
$ws.Range("N138").Value2 = -21445.5716

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value2 = 99.5
$ws.Range("I4").Value2 = 99.5
$ws.Range("J4").Value2 = 0
$ws.Range("K4").Value2 = 99.5
$ws.Range("L4").Value2 = 0
$ws.Range("M4").ClearContents()
$ws.Range("N4").Value2 = 16.5
$ws.Range("H6").Value2 = 21000000
$ws.Range("I6").Value2 = 19166666
$ws.Range("K6").Value2 = 19166666
$ws.Range("M6").Value2 = -19166493
$ws.Range("H37").Value2 = 99518.5
$ws.Range("J37").Value2 = 99518.5
$ws.Range("L37").Value2 = 99518.5
$ws.Range("N37").Value2 = -100064.5
$ws.Range("H41").Value2 = 1410.25
$ws.Range("J41").Value2 = 4000
$ws.Range("L41").Value2 = 4000
$ws.Range("N41").Value2 = -4828
$ws.Range("H122").Value2 = 3019.3845
$ws.Range("I122").Value2 = 2973.2222
$ws.Range("J122").Value2 = 3123.25
$ws.Range("K122").Value2 = 8919.6666
$ws.Range("L122").Value2 = 9369.75
$ws.Range("M122").Value2 = -6469.6666
$ws.Range("N122").Value2 = -14269.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H11").Value2 = 808.3333
$ws.Range("I11").Value2 = 222.5
$ws.Range("K11").Value2 = 222.5
$ws.Range("M11").Value2 = -82.5
$ws.Range("H134").Value2 = 1937.0952
$ws.Range("I134").Value2 = 1382.5333
$ws.Range("J134").Value2 = 3323.5
$ws.Range("K134").Value2 = 4147.5999
$ws.Range("L134").Value2 = 9970.5
$ws.Range("M134").Value2 = -1612.5999
$ws.Range("N134").Value2 = -15040.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value2 = 1245.3077
$ws.Range("I122").Value2 = 1160.1428
$ws.Range("J122").Value2 = 1344.6666
$ws.Range("K122").Value2 = 3480.4284
$ws.Range("L122").Value2 = 4033.9998
$ws.Range("M122").Value2 = -1030.4284
$ws.Range("N122").Value2 = -8933.9998
$ws.Range("H132").Value2 = 3259.2666
$ws.Range("I132").Value2 = 2525.55
$ws.Range("K132").Value2 = 7576.650000000001
$ws.Range("M132").Value2 = -5046.650000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value2 = 32.52174
$ws.Range("I2").Value2 = 32.7
$ws.Range("K2").Value2 = 196.2
$ws.Range("M2").Value2 = -83.20000000000002
$ws.Range("H7").Value2 = 56.25
$ws.Range("J7").Value2 = 100
$ws.Range("L7").Value2 = 300
$ws.Range("N7").Value2 = -524
$ws.Range("H24").Value2 = 400
$ws.Range("J24").Value2 = 400
$ws.Range("L24").Value2 = 1200
$ws.Range("N24").Value2 = -1660
$ws.Range("H34").Value2 = 2248.4614
$ws.Range("J34").Value2 = 2352.3333
$ws.Range("L34").Value2 = 7056.999899999999
$ws.Range("N34").Value2 = -7224.999899999999
$ws.Range("H39").Value2 = 5811.1113
$ws.Range("J39").Value2 = 6437.5
$ws.Range("L39").Value2 = 19312.5
$ws.Range("N39").Value2 = -19900.5
$ws.Range("H55").Value2 = 5012.727
$ws.Range("J55").Value2 = 5012.727
$ws.Range("L55").Value2 = 15038.181
$ws.Range("N55").Value2 = -15392.181
$ws.Range("H92").Value2 = 7375.75
$ws.Range("J92").Value2 = 5251.5
$ws.Range("L92").Value2 = 15754.5
$ws.Range("N92").Value2 = -18250.5
$ws.Range("H117").Value2 = 300
$ws.Range("J117").Value2 = 300
$ws.Range("L117").Value2 = 900
$ws.Range("N117").Value2 = -7784
$ws.Range("H122").Value2 = 1128.3334
$ws.Range("I122").Value2 = 1099
$ws.Range("K122").Value2 = 9891
$ws.Range("M122").Value2 = -7441

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value2 = 3558.5625
$ws.Range("I122").Value2 = 2360
$ws.Range("J122").Value2 = 6195.4
$ws.Range("K122").Value2 = 7080
$ws.Range("L122").Value2 = 18586.2
$ws.Range("M122").Value2 = -4630
$ws.Range("N122").Value2 = -23486.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value2 = 6108.25
$ws.Range("I7").Value2 = 4493.7144
$ws.Range("K7").Value2 = 4493.7144
$ws.Range("M7").Value2 = -4381.7144
$ws.Range("H31").Value2 = 573.8461
$ws.Range("I31").Value2 = 296.1111
$ws.Range("J31").Value2 = 1198.75
$ws.Range("K31").Value2 = 296.1111
$ws.Range("L31").Value2 = 1198.75
$ws.Range("M31").Value2 = -48.11110000000002
$ws.Range("N31").Value2 = -1694.75
$ws.Range("H46").Value2 = 4268.6562
$ws.Range("J46").Value2 = 4584.476
$ws.Range("L46").Value2 = 4584.476
$ws.Range("N46").Value2 = -4960.476
$ws.Range("H55").Value2 = 1116.75
$ws.Range("I55").Value2 = 1067.4445
$ws.Range("J55").Value2 = 1180.1428
$ws.Range("K55").Value2 = 1067.4445
$ws.Range("L55").Value2 = 1180.1428
$ws.Range("M55").Value2 = -894.4445000000001
$ws.Range("N55").Value2 = -1526.1428
$ws.Range("H126").Value2 = 6108.25
$ws.Range("I126").Value2 = 4493.7144
$ws.Range("K126").Value2 = 13481.1432
$ws.Range("M126").Value2 = -11011.1432

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value2 = 1920.6666
$ws.Range("I122").Value2 = 1504.8
$ws.Range("J122").Value2 = 4000
$ws.Range("K122").Value2 = 4514.4
$ws.Range("L122").Value2 = 12000
$ws.Range("M122").Value2 = -2064.4
$ws.Range("N122").Value2 = -16900
$ws.Range("H126").Value2 = 6178.9
$ws.Range("I126").Value2 = 3447.25
$ws.Range("K126").Value2 = 10341.75
$ws.Range("M126").Value2 = -7871.75
